$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (the "Powerp Odoo" asset) was missing its purchase_date. Give F4
# the same look (number format / style) as the other purchase_date cells
# in column F (F2, F3) by copying the format over, then write the date
# text value that goes with it.
$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("F4").Value = "<2-10-01"

# Reflect the new active selection recorded in the saved sheet.
$ws.Range("F4").Select()
